$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (border/bold/alignment) from H1 into the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF)
$iValues = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1;
    19 = 1; 20 = 1; 21 = 5; 22 = 9; 23 = 4; 24 = 6; 25 = 4
}
$jValues = @{
    2 = 3; 3 = 7; 4 = 6; 5 = 6; 6 = 6; 7 = 7; 8 = 5; 9 = 4; 10 = 4;
    11 = 6; 12 = 5; 13 = 6; 14 = 4; 15 = 7; 16 = 6; 17 = 6; 18 = 6;
    19 = 6; 20 = 5; 21 = 9; 22 = 9; 23 = 5; 24 = 6; 25 = 5
}

for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
